# The commit reshuffles the SharePoint-managed custom XML parts that ship
# inside the package: the "document management properties" part
# (root element p:properties, schema http://schemas.microsoft.com/office/2006/metadata/properties)
# and the "content type schema" part
# (root element ct:contentTypeSchema, schema http://schemas.microsoft.com/office/2006/metadata/contentType)
# trade physical slots (customXml/item2.xml <-> customXml/item3.xml, and their
# matching customXml/itemPropsN.xml companions). The actual XML payload of each
# part is untouched - only which numbered part holds which payload changes.
#
# Do this through the real object model: pull each part's XML by its distinct
# namespace, delete both parts, then re-add them in the opposite order so the
# package numbers them the other way around (Office assigns the item*.xml /
# itemProps*.xml slots in the order parts are (re)created).

$p = $ppt.ActivePresentation
$parts = $p.CustomXMLParts

$propertiesNs = "http://schemas.microsoft.com/office/2006/metadata/properties"
$contentTypeNs = "http://schemas.microsoft.com/office/2006/metadata/contentType"

$propertiesXml = $null
$contentTypeXml = $null
$propertiesPart = $null
$contentTypePart = $null

for ($i = 1; $i -le $parts.Count; $i++) {
    $part = $parts.Item($i)
    $ns = $part.NamespaceURI
    if ($ns -eq $propertiesNs) {
        $propertiesPart = $part
        $propertiesXml = $part.XML
    } elseif ($ns -eq $contentTypeNs) {
        $contentTypePart = $part
        $contentTypeXml = $part.XML
    }
}

if (($propertiesPart -ne $null) -and ($contentTypePart -ne $null)) {
    # Remove both existing parts ...
    $propertiesPart.Delete()
    $contentTypePart.Delete()

    # ... and recreate them swapped, so the content-type schema now lands in
    # the slot the properties part used to occupy, and vice versa.
    $parts.Add($contentTypeXml, $contentTypeNs)
    $parts.Add($propertiesXml, $propertiesNs)
}
